$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price (D) column cells we touch so Excel does not
# auto-convert values such as "1.001" or "26.442.26" into numbers/dates.
$priceCells = @("D2","D3","D5","D7","D9","D10","D11","D12","D13","D14","D15","D17","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.442.26'
$ws.Range('E2').Value = '  +3.44%  '
$ws.Range('D3').Value = '1.729.65'
$ws.Range('E3').Value = '  +3.90%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').Value = '243.41'
$ws.Range('E5').Value = '  +2.76%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '0.4792'
$ws.Range('E7').Value = '  +3.79%  '
$ws.Range('E8').Value = '  +3.42%  '
$ws.Range('D9').Value = '0.06221'
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('D10').Value = '1.731.65'
$ws.Range('E10').Value = '  +4.06%  '
$ws.Range('D11').Value = '0.07122'
$ws.Range('E11').Value = '  +2.82%  '
$ws.Range('D12').Value = '15.67'
$ws.Range('E12').Value = '  +5.88%  '
$ws.Range('D13').Value = '0.6146'
$ws.Range('E13').Value = '  +7.37%  '
$ws.Range('D14').Value = '4.548'
$ws.Range('E14').Value = '  +5.06%  '
$ws.Range('D15').Value = '76.90'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('E16').Value = '  +0.09%  '
$ws.Range('D17').Value = '26.468.47'
$ws.Range('E17').Value = '  +3.51%  '
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = '0.000006908'
$ws.Range('E19').Value = '  +3.49%  '
$ws.Range('D20').Value = '11.69'
$ws.Range('E20').Value = '  +2.99%  '
$ws.Range('D21').Value = '1.956.13'
$ws.Range('E21').Value = '  +4.23%  '
$ws.Range('D22').Value = '4.569'
$ws.Range('E22').Value = '  +3.32%  '
$ws.Range('D23').Value = '8.891'
$ws.Range('E23').Value = '  +3.52%  '
$ws.Range('D24').Value = '5.319'
$ws.Range('E24').Value = '  +2.00%  '
$ws.Range('D25').Value = '136.40'
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('D26').Value = '15.33'
$ws.Range('E26').Value = '  +2.79%  '
$ws.Range('D27').Value = '1.788'
$ws.Range('E27').Value = '  +4.33%  '
$ws.Range('D28').Value = '1.403'
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('D29').Value = '106.54'
$ws.Range('E29').Value = '  +2.52%  '
$ws.Range('D30').Value = '3.975'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('D31').Value = '0.08014'
$ws.Range('E31').Value = '  +4.88%  '
$ws.Range('D32').Value = '3.713'
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('D33').Value = '0.04534'
$ws.Range('E33').Value = '  +4.62%  '
$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').Value = '1.001'
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.618'
$ws.Range('E35').Value = '  +0.50%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.6345'
$ws.Range('E36').Value = '  +4.81%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '0.9866'
$ws.Range('E37').Value = '  +5.21%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '0.9317'
$ws.Range('E38').Value = '  +0.42%  '
$ws.Range('B39').Value = 'Quant'
$ws.Range('C39').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D39').Value = '109.59'
$ws.Range('E39').Value = '  +2.61%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '1.974'
$ws.Range('E40').Value = '  +8.08%  '
$ws.Range('B41').Value = 'MXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D41').Value = '2.394'
$ws.Range('E41').Value = '  -1.56%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '1.006'
$ws.Range('E42').Value = '  +0.68%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.01504'
$ws.Range('E43').Value = '  +4.23%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '5.642'
$ws.Range('E44').Value = '  +11.75%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '0.3896'
$ws.Range('E45').Value = '  +5.29%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '6.931'
$ws.Range('E46').Value = '  +14.26%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1189'
$ws.Range('E47').Value = '  +7.23%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.05329'
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('D49').Value = '30.80'
$ws.Range('E49').Value = '  -1.33%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '7.852'
$ws.Range('E50').Value = '  +3.78%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.271'
